# Update panels datasheet (Sheet1) and data entry spreadsheet (Sheet1 (2))
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet1 (2)")

# ---------------------------------------------------------------------
# 1) Sheet1: move the "species list" sub-header (Container Weight /
#    Container + Wet Weight / Wet Weight) up from row 28 to row 19, and
#    widen the 3D-photos entry rows (20-25) to cover columns C:E too.
# ---------------------------------------------------------------------

# Grab the formatting (border/alignment) that currently lives on C28:E28
# and stamp it onto C19:E19 before we touch the text in row 28.
$ws1.Range("C28:E28").Copy()
$ws1.Range("C19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("C19").Value = "Container Weight (grams)"
$ws1.Range("D19").Value = "Container + Wet Weight (grams)"
$ws1.Range("E19").Value = "Wet Weight (grams)"

# Row 19 grows taller to fit the wrapped header text.
$ws1.Rows.Item(19).RowHeight = 30.75

# Rows 20-25 (the "1)".."6)" panel rows) pick up blank, centre-bordered
# cells in C:E matching the existing blank-cell style used elsewhere
# (e.g. B1), and get a touch more height.
$ws1.Range("B1").Copy()
$ws1.Range("C20:E25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Rows.Item(20).RowHeight = 27.75
$ws1.Rows.Item(21).RowHeight = 27.75
$ws1.Rows.Item(22).RowHeight = 27.75
$ws1.Rows.Item(23).RowHeight = 27.75
$ws1.Rows.Item(24).RowHeight = 27.75
$ws1.Rows.Item(25).RowHeight = 27.75

# ---------------------------------------------------------------------
# 2) Row 28 keeps its border formatting but loses its header text (the
#    text moved up to row 19), and its row height shrinks back down
#    since it no longer needs to wrap.
# ---------------------------------------------------------------------
$ws1.Range("C28:E28").ClearContents()
$ws1.Rows.Item(28).RowHeight = 22.5

# Rows 29-31 (the eDNA entry rows) switch their C:E formatting from the
# bordered "B1-style" cells to the borderless vertical-center style
# already used at B17.
$ws1.Range("B17").Copy()
$ws1.Range("C29:E31").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New row 32: a thin top-border divider cell in B32, and a plain
# borderless vertical-center cell in C32 (matching B17's style).
$ws1.Range("B17").Copy()
$ws1.Range("C32").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("B32").Borders.Item(8).LineStyle = 1
$ws1.Range("B32").VerticalAlignment = -4108

# ---------------------------------------------------------------------
# 3) Page setup: Sheet1 print scale drops from 94% to 90%.
# ---------------------------------------------------------------------
$ws1.PageSetup.Zoom = $false
$ws1.PageSetup.FitToPagesWide = $false
$ws1.PageSetup.FitToPagesTall = $false
$ws1.PageSetup.Zoom = 90

# ---------------------------------------------------------------------
# 4) Selection / active sheet bookkeeping: Sheet1's frozen topLeftCell
#    scroll position is cleared and its selection moves to A20; Sheet1
#    (2) becomes the active sheet with its selection also at A20.
# ---------------------------------------------------------------------
$ws1.Range("A20").Select()
$ws2.Activate()
$ws2.Range("A20").Select()
